$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (bold, centered, top-aligned, thin border) used by column A
# "Predicted"/"Baseline" labels (same style as existing A2:A17) onto the new rows.
$ws.Range("A17").Copy()
$ws.Range("A18:A33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range('A18').Value = 'Predicted'
$ws.Range('B18').Value = '
"heritage tourism" OR "nature tourism" OR "tourist attractions" OR "hospitality industry" OR "economic impact of tourism" OR "rural tourism" OR "urban tourism" OR "ecotourism" OR "tourism infrastructure" OR "community-based tourism" OR "tourism investment" OR "medical tourism" OR "digital tourism" OR "tourism statistics" OR "responsible tourism" OR "seasonal tourism" OR "adventure tourism" OR (tourism growth nexus) OR 
(("tourism development" OR "tourism management" OR "tourism marketing" OR "cultural tourism" OR "travel technology" OR "sustainable tourism" OR "destination management" OR "tourism policy" OR "business tourism" OR "global tourism" OR "travel trends" OR "tourist behavior") AND (Nexus))
'
$ws.Range('C18').Value = 0.08
$ws.Range('D18').Value = 0
$ws.Range('E18').Value = 0.001
$ws.Range('F18').Value = 0.413
$ws.Range('G18').Value = 0.095
$ws.Rows.Item(18).AutoFit()

$ws.Range('A19').Value = 'Baseline'
$ws.Range('B19').Value = 'Tourism Growth Nexus'
$ws.Range('C19').Value = 0.08
$ws.Range('D19').Value = 0.019
$ws.Range('E19').Value = 0.049
$ws.Range('F19').Value = 0.421
$ws.Range('G19').Value = 0.095
$ws.Rows.Item(19).AutoFit()

$ws.Range('A20').Value = 'Predicted'
$ws.Range('B20').Value = '
"exergy analysis" OR "sustainable biofuel" OR "lifecycle assessment" OR "carbon neutrality" OR "algal biofuels" OR "food vs fuel debate" OR 
(("waste-to-energy" OR "feedstock" OR "biodiesel" OR "economic viability" OR "bioenergy" OR "agricultural residues" OR "fossil fuel alternatives" OR "circular economy" OR "bioethanol" OR "energy independence") AND (Biofuel))
'
$ws.Range('C20').Value = 0.36
$ws.Range('D20').Value = 0.454
$ws.Range('E20').Value = 0.375
$ws.Range('F20').Value = 0.413
$ws.Range('G20').Value = 0.369
$ws.Rows.Item(20).AutoFit()

$ws.Range('A21').Value = 'Baseline'
$ws.Range('B21').Value = 'Sustainable Biofuel Economy'
$ws.Range('C21').Value = 0.02
$ws.Range('D21').Value = 0.151
$ws.Range('E21').Value = 0.024
$ws.Range('F21').Value = 0.467
$ws.Range('G21').Value = 0.025
$ws.Rows.Item(21).AutoFit()

$ws.Range('A22').Value = 'Predicted'
$ws.Range('B22').Value = '
"perovskite solar cells" OR "perovskite composition" OR "moisture resistance" OR "halide perovskites" OR "bandgap engineering" OR "photoelectric conversion" OR "mapbi3" OR "fapbi3" OR 
(("efficiency" OR "stability" OR "degradation" OR "thermal stability" OR "nucleation" OR "doping effects" OR "fabrication techniques" OR "molecular dynamics simulations" OR "structural stability" OR "defect states" OR "scalability" OR "environmental stability" OR "lifespan" OR "energy conversion efficiency" OR "crystallization process" OR "interface stability" OR "photovoltaic performance" OR "solar energy conversion" OR "optoelectronic properties" OR "thin-film technology" OR "surface passivation" OR "charge carrier dynamics") AND (perovskite AND Solar))
'
$ws.Range('C22').Value = 1
$ws.Range('D22').Value = 0.524
$ws.Range('E22').Value = 0.846
$ws.Range('F22').Value = 0.461
$ws.Range('G22').Value = 0.8110000000000001
$ws.Rows.Item(22).AutoFit()

$ws.Range('A23').Value = 'Baseline'
$ws.Range('B23').Value = 'Perovskite Solar Cells Stability'
$ws.Range('C23').Value = 0.897
$ws.Range('D23').Value = 0.8149999999999999
$ws.Range('E23').Value = 0.879
$ws.Range('F23').Value = 0.481
$ws.Range('G23').Value = 0.764
$ws.Rows.Item(23).AutoFit()

$ws.Range('A24').Value = 'Predicted'
$ws.Range('B24').Value = '
"nanocarriers" OR "nanoparticle characterization" OR "translational medicine" OR "camptothecin" OR "ribonucleotide reductase inhibitors" OR (nanopharmaceuticals) OR (nanonutraceuticals) OR 
(("pharmacodynamics" OR "personalized medicine") AND (nanoparticles))
'
$ws.Range('C24').Value = 0.04
$ws.Range('D24').Value = 0.007
$ws.Range('E24').Value = 0.021
$ws.Range('F24').Value = 0.534
$ws.Range('G24').Value = 0.049
$ws.Rows.Item(24).AutoFit()

$ws.Range('A25').Value = 'Baseline'
$ws.Range('B25').Value = 'Nanopharmaceuticals OR Nanonutraceuticals'
$ws.Range('C25').Value = 0
$ws.Range('D25').Value = 0
$ws.Range('E25').Value = 0
$ws.Range('F25').Value = 0.583
$ws.Range('G25').Value = 0
$ws.Rows.Item(25).AutoFit()

$ws.Range('A26').Value = 'Predicted'
$ws.Range('B26').Value = '
"climate action" OR "logistics efficiency" OR "circular economy" OR "sustainable logistics" OR "sustainable supply chain" OR "supply chain sustainability" OR "green building design" OR "inventory optimization" OR "green transportation" OR "eco-efficient processes" OR "carbon footprint reduction" OR "green certifications" OR "sustainable procurement" OR "warehouse automation" OR "eco-friendly packaging" OR "leed certification" OR "smart warehousing" OR "last mile delivery solutions" OR (green warehousing) OR "cold chain sustainability" OR 
(("energy efficiency" OR "renewable energy" OR "emission reduction" OR "environmental impact assessment" OR "energy management systems" OR "waste management" OR "sustainable operations" OR "transport emissions" OR "green technologies" OR "resource optimization") AND (Warehousing))
'
$ws.Range('C26').Value = 0.237
$ws.Range('D26').Value = 0.015
$ws.Range('E26').Value = 0.061
$ws.Range('F26').Value = 0.432
$ws.Range('G26').Value = 0.26
$ws.Rows.Item(26).AutoFit()

$ws.Range('A27').Value = 'Baseline'
$ws.Range('B27').Value = 'Green Warehousing'
$ws.Range('C27').Value = 0.158
$ws.Range('D27').Value = 0.422
$ws.Range('E27').Value = 0.181
$ws.Range('F27').Value = 0.412
$ws.Range('G27').Value = 0.18
$ws.Rows.Item(27).AutoFit()

$ws.Range('A28').Value = 'Predicted'
$ws.Range('B28').Value = '
"edge data processing" OR "distributed machine learning" OR "hardware acceleration" OR "deep learning inference" OR "edge infrastructure" OR "federated learning" OR "context-aware computing" OR "resource-constrained environments" OR "sensor data analytics" OR "smart edge applications" OR "edge-to-cloud architecture" OR "on-device ai" OR "latency-sensitive applications" OR "edge analytics" OR "ai model deployment" OR "edge ai devices" OR (ai on edge devices) OR "energy-efficient ai" OR "cognitive edge computing" OR "privacy-preserving ai" OR "ai inference engines" OR "scalability in edge ai" OR "ai multi-tenancy" OR 
(("autonomous systems" OR "real-time data processing" OR "model compression" OR "internet of things (iot)" OR "low latency" OR "edge computing" OR "neural network optimization" OR "heterogeneous computing") AND (Edge Devices))
'
$ws.Range('C28').Value = 0.389
$ws.Range('D28').Value = 0.161
$ws.Range('E28').Value = 0.303
$ws.Range('F28').Value = 0.522
$ws.Range('G28').Value = 0.41
$ws.Rows.Item(28).AutoFit()

$ws.Range('A29').Value = 'Baseline'
$ws.Range('B29').Value = 'AI on Edge Devices'
$ws.Range('C29').Value = 0.111
$ws.Range('D29').Value = 0.116
$ws.Range('E29').Value = 0.112
$ws.Range('F29').Value = 0.483
$ws.Range('G29').Value = 0.131
$ws.Rows.Item(29).AutoFit()

$ws.Range('A30').Value = 'Predicted'
$ws.Range('B30').Value = '
"smart healthcare" OR "mhealth" OR (internet of things in healthcare) OR "medical device integration" OR "internet of medical things" OR "smart hospitals" OR "ai in healthcare" OR "real-time health tracking" OR "healthcare interoperability" OR 
(("chronic disease management" OR "health monitoring systems" OR "population health management" OR "wearable devices" OR "connected health" OR "digital health" OR "interoperability" OR "electronic health records" OR "patient engagement" OR "telemedicine" OR "virtual health" OR "personalized medicine" OR "remote patient monitoring" OR "health information exchange" OR "predictive analytics" OR "secure health data" OR "digital therapeutics" OR "health data analytics" OR "healthcare automation") AND (IoT))
'
$ws.Range('C30').Value = 0.517
$ws.Range('D30').Value = 0.157
$ws.Range('E30').Value = 0.354
$ws.Range('F30').Value = 0.676
$ws.Range('G30').Value = 0.543
$ws.Rows.Item(30).AutoFit()

$ws.Range('A31').Value = 'Baseline'
$ws.Range('B31').Value = 'Internet of Things in Healthcare'
$ws.Range('C31').Value = 0.345
$ws.Range('D31').Value = 0.326
$ws.Range('E31').Value = 0.341
$ws.Range('F31').Value = 0.596
$ws.Range('G31').Value = 0.377
$ws.Rows.Item(31).AutoFit()

$ws.Range('A32').Value = 'Predicted'
$ws.Range('B32').Value = '
"bioenergy" OR "decarbonization" OR (energy growth nexus)
'
$ws.Range('C32').Value = 0.037
$ws.Range('D32').Value = 0.032
$ws.Range('E32').Value = 0.036
$ws.Range('F32').Value = 0.588
$ws.Range('G32').Value = 0.046
$ws.Rows.Item(32).AutoFit()

$ws.Range('A33').Value = 'Baseline'
$ws.Range('B33').Value = 'Energy Growth Nexus'
$ws.Range('C33').Value = 0.037
$ws.Range('D33').Value = 0.421
$ws.Range('E33').Value = 0.045
$ws.Range('F33').Value = 0.523
$ws.Range('G33').Value = 0.045
$ws.Rows.Item(33).AutoFit()

$ws.Range("A1").Select()